$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Oven-dry total biomass weights: row, B (OD_Total Biomass_Bag_Tie_g), C (Bag_Tie_g)
$data = @(
    @(2, 1574.4, 120),
    @(3, 1210.2, 121),
    @(4, 1359.7, 124),
    @(5, 1302.8, 119),
    @(6, 1220.8, 120),
    @(7, 1165.5999999999999, 122),
    @(8, 1647.4, 124),
    @(9, 1353.8, 121),
    @(10, 1430.6, 179),
    @(11, 1335.5, 125),
    @(12, 1360.5, 122),
    @(13, 1376.1, 125),
    @(14, 1531.4, 123),
    @(15, 1469.9, 126),
    @(16, 1749.5, 127),
    @(17, 1315.4, 119),
    @(18, 1366.4, 118),
    @(19, 1450.9, 119)
)

foreach ($row in $data) {
    $r = $row[0]
    $b = $row[1]
    $c = $row[2]
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
}

# D2 gets its own formula, D3:D19 are filled as one block (mirrors the
# "fill down" shared-formula pattern Excel produces for a dragged formula).
$ws.Range("D2").Formula = "=B2-C2"
$ws.Range("D3:D19").Formula = "=B3-C3"

$ws.Range("D20").Select() | Out-Null
